$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force target cells to Text format so numeric-looking strings (e.g. "1.004")
# and values with trailing zeros are preserved exactly as text, matching the
# original inlineStr cell content.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('B13').NumberFormat = '@'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('B14').NumberFormat = '@'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('E51').NumberFormat = '@'

$ws.Range('D2').Value = '25.857.12'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').Value = '1.636.53'
$ws.Range('E3').Value = '  +0.71%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Value = '215.15'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').Value = '0.5088'
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').Value = '1.004'
$ws.Range('E7').Value = '  +0.26%  '
$ws.Range('D8').Value = '0.2583'
$ws.Range('E8').Value = '  +1.01%  '
$ws.Range('D9').Value = '0.06427'
$ws.Range('E9').Value = '  +1.78%  '
$ws.Range('E10').Value = '  +5.10%  '
$ws.Range('D11').Value = '0.07798'
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('D12').Value = '4.270'
$ws.Range('E12').Value = '  +1.16%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.643.94'
$ws.Range('E13').Value = '  +1.10%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '1.861.89'
$ws.Range('E14').Value = '  +0.80%  '
$ws.Range('D15').Value = '0.5601'
$ws.Range('E15').Value = '  +1.50%  '
$ws.Range('D16').Value = '0.0₅7662'
$ws.Range('E16').Value = '  +2.31%  '
$ws.Range('D17').Value = '63.24'
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('D18').Value = '25.867.63'
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('D19').Value = '1.004'
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('D20').Value = '193.42'
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').Value = '4.385'
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('D22').Value = '9.946'
$ws.Range('E22').Value = '  +2.06%  '
$ws.Range('D23').Value = '6.156'
$ws.Range('E23').Value = '  +2.72%  '
$ws.Range('D24').Value = '1.004'
$ws.Range('E24').Value = '  +0.29%  '
$ws.Range('D25').Value = '1.789'
$ws.Range('E25').Value = '  -4.26%  '
$ws.Range('D26').Value = '138.67'
$ws.Range('D27').Value = '0.1230'
$ws.Range('E27').Value = '  -0.75%  '
$ws.Range('D28').Value = '6.854'
$ws.Range('E28').Value = '  +2.36%  '
$ws.Range('E29').Value = '  +0.85%  '
$ws.Range('D30').Value = '1.240'
$ws.Range('E30').Value = '  +0.22%  '
$ws.Range('D31').Value = '0.04951'
$ws.Range('E31').Value = '  +2.18%  '
$ws.Range('D32').Value = '3.298'
$ws.Range('E32').Value = '  +2.02%  '
$ws.Range('D33').Value = '3.250'
$ws.Range('E33').Value = '  +3.02%  '
$ws.Range('D34').Value = '1.567'
$ws.Range('E34').Value = '  +2.21%  '
$ws.Range('D35').Value = '2.388'
$ws.Range('E35').Value = '  +1.03%  '
$ws.Range('D36').Value = '0.9037'
$ws.Range('E36').Value = '  +1.31%  '
$ws.Range('D37').Value = '2.580'
$ws.Range('E37').Value = '  +1.70%  '
$ws.Range('D38').Value = '0.5557'
$ws.Range('E38').Value = '  +1.22%  '
$ws.Range('D39').Value = '1.134.12'
$ws.Range('E39').Value = '  +1.92%  '
$ws.Range('D40').Value = '0.01570'
$ws.Range('E40').Value = '  +1.59%  '
$ws.Range('D41').Value = '0.9973'
$ws.Range('E41').Value = '  -0.34%  '
$ws.Range('D42').Value = '99.25'
$ws.Range('E42').Value = '  +2.32%  '
$ws.Range('D43').Value = '5.464'
$ws.Range('E43').Value = '  -0.86%  '
$ws.Range('D44').Value = '0.8018'
$ws.Range('E44').Value = '  +0.82%  '
$ws.Range('D45').Value = '0.0₈111'
$ws.Range('E45').Value = '  -3.36%  '
$ws.Range('D46').Value = '55.48'
$ws.Range('E46').Value = '  +1.85%  '
$ws.Range('D47').Value = '0.4261'
$ws.Range('E47').Value = '  -3.67%  '
$ws.Range('D48').Value = '7.779'
$ws.Range('E48').Value = '  +3.05%  '
$ws.Range('D49').Value = '0.05070'
$ws.Range('E49').Value = '  -1.10%  '
$ws.Range('D50').Value = '0.9998'
$ws.Range('E50').Value = '  +0.39%  '
$ws.Range('D51').Value = '1.003'
$ws.Range('E51').Value = '  +0.27%  '
